# Refresh the cryptocurrency snapshot table (columns B:E, rows 2-51) on Sheet1
# to match the latest coinranking.com figures, per the automated GitHub Actions update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.550.70'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.39%  '
# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.417.94'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.29%  '
# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.78'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.68%  '
# Row 6: Solana
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -3.21%  '
# Row 7: XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.620'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +7.45%  '
# Row 8: USDC
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.06%  '
# Row 9: LidoStakedEther
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.425.41'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.08%  '
# Row 10: Toncoin
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.23%  '
# Row 11: Dogecoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.123'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.74%  '
# Row 12: Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.442'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.51%  '
# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.008.69'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.18%  '
# Row 14: TRON
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.19%  '
# Row 15: ShibaInu
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.45%  '
# Row 16: Avalanche
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.96'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.93%  '
# Row 17: WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.571.90'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.30%  '
# Row 18: WrappedEther
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.433.12'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.89%  '
# Row 19: Polkadot
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.35'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.33%  '
# Row 20: Chainlink
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.51%  '
# Row 21: BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '375.42'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.73%  '
# Row 22: Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.01'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.35%  '
# Row 23: Polygon
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.03%  '
# Row 24: Dai (was Litecoin row)
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.27%  '
# Row 25: Litecoin (was Dai row)
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.48'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.75%  '
# Row 26: PEPE
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.42%  '
# Row 27: InternetComputer(DFINITY)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.28'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +8.30%  '
# Row 28: Kaspa
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.54%  '
# Row 29: Binance-PegBSC-USD
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.07%  '
# Row 30: Fetch.AI
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.00%  '
# Row 31: NEARProtocol
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.14%  '
# Row 32: PancakeSwap
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.03'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.87%  '
# Row 33: EthereumClassic
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.12'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.19%  '
# Row 34: Aptos
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.23'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.78%  '
# Row 35: ImmutableX
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +6.26%  '
# Row 36: Monero
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '160.49'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.94%  '
# Row 37: Stacks
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.90'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.77%  '
# Row 38: RenderToken (was Hedera row)
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.98'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +6.65%  '
# Row 39: Hedera (was RenderToken row)
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0762'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.39%  '
# Row 40: EnergySwap
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '26.63'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.88%  '
# Row 41: Maker
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.864.32'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.35%  '
# Row 42: Filecoin
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.62'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.20%  '
# Row 43: OKB
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.92'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.19%  '
# Row 44: InjectiveProtocol
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '26.43'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +9.63%  '
# Row 45: VeChain
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0315'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.68%  '
# Row 46: Mantle
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.770'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.83%  '
# Row 47: Bittensor
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '319.74'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +7.34%  '
# Row 48: ONDO
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.72%  '
# Row 50: Cosmos
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.58'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.96%  '
# Row 51: SuiNetwork
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.856'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.93%  '
